$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 50
$ws.Cells.Item($row, 1).Value = "G3"
$ws.Cells.Item($row, 2).Value = "Eat Healthy"
$ws.Cells.Item($row, 3).Value = 45907
$ws.Cells.Item($row, 3).NumberFormat = $ws.Cells.Item($row - 1, 3).NumberFormat
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
